# Apply the Translation-sheet reshuffle described by the diff.
# (Commit message: "Issue was solved by changing lcd24bpp to lcd32bpp" --
#  the underlying texts.xlsx is an auto-generated TouchGFX asset; the
#  visible effect in this workbook is a re-ordering of the per-text
#  rows (TEXT ID / ALIGNMENT / GB translation) on the "Translation" sheet,
#  rows 4..14, columns B/D/F. Columns C and E are untouched.)
#
# NOTE: `.Value` (read) is unreliable in this host -- use `.Value2` to read.
# `.Value` (write) auto-converts a numeric-looking string ("1", "2", ...)
# to a real number, so for those destination cells we briefly force a Text
# number format before writing, then restore the style afterwards so the
# saved cell carries no extra `s=` attribute.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Snapshot the current (pre-edit) B/D/F text for rows 4..14.
$colB = @{}
$colD = @{}
$colF = @{}
for ($r = 4; $r -le 14; $r++) {
    $colB[$r] = $ws.Cells.Item($r, 2).Value2
    $colD[$r] = $ws.Cells.Item($r, 4).Value2
    $colF[$r] = $ws.Cells.Item($r, 6).Value2
}

# New row -> old row it should take its (B, D, F) content from.
$sourceRow = @{
    4  = 6
    5  = 7
    6  = 8
    7  = 9
    8  = 10
    9  = 5
    10 = 13
    11 = 14
    12 = 4
    13 = 12
    14 = 11
}

# Rows whose new F value is purely numeric text ("1".."4") and therefore
# needs to be pinned to Text format before assignment so it doesn't turn
# into a real number.
$numericFRows = @(5, 6, 7, 8)

foreach ($r in $numericFRows) {
    $ws.Cells.Item($r, 6).NumberFormat = "@"
}

foreach ($r in $sourceRow.Keys) {
    $src = $sourceRow[$r]
    $ws.Cells.Item($r, 2).Value = $colB[$src]
    $ws.Cells.Item($r, 4).Value = $colD[$src]
    $ws.Cells.Item($r, 6).Value = $colF[$src]
}

foreach ($r in $numericFRows) {
    $ws.Cells.Item($r, 6).Style = "Normal"
}
